# [Update] 타격 VFX 추가
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Hit_vfx prefab path (B15) to point at the new "Holy hit" prefab.
$ws.Range("B15").Value = "Assets/Resource/VFX/PlayerVFX/Hit_vfx/Hit_vfx/Holy hit.prefab"

# Move the active selection to D13, matching the saved workbook view state.
$ws.Range("D13").Select()
